# Apply the Alvearie FHIR IG CodeSystem-identifier-type.xlsx update:
#  - Version: 6.1.0 -> 6.1.1
#  - Date: 2022-05-31T20:10:14+00:00 -> 2022-06-06T15:56:40+00:00
#  - Count: 6 -> 7
#  - Concepts sheet: add a new "TKN" / "Token identifier" row

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet updates -------------------------------------------------

# Version (row 3)
$meta.Range("B3").Value = "6.1.1"

# Date (row 8)
$meta.Range("B8").Value = "2022-06-06T15:56:40+00:00"

# Count (row 17) - value is numeric-looking text ("7"), so force it to stay
# text the same way the original "6" was stored (as a shared string), using a
# scratch cell so the text number format doesn't get attached to the real
# cell's style.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "7"
$meta.Range("Z1").Copy()
$meta.Range("B17").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

# --- Concepts sheet: append new concept row ---------------------------------

$newRow = 8

# Level value ("1") is also numeric-looking text, use the same scratch-cell
# trick to keep it a text/shared-string value like the rest of the column.
$concepts.Range("Z1").NumberFormat = "@"
$concepts.Range("Z1").Value = "1"
$concepts.Range("Z1").Copy()
$concepts.Range("A" + $newRow).PasteSpecial(-4163)
$concepts.Range("Z1").Clear()

$concepts.Range("B" + $newRow).Value = "TKN"
$concepts.Range("C" + $newRow).Value = "Token identifier"
$concepts.Range("D" + $newRow).Value = "Data tokenization service token"

# Copy the row-7 formatting onto the new row so it matches the rest of the
# table (border/fill/alignment style).
$concepts.Range("A7:D7").Copy()
$concepts.Range("A8:D8").PasteSpecial(-4122)
